$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = "'05/08/2025"
$ws.Range("A23").Style = "Normal"
$ws.Range("B23").Value = "Cienciano"
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = "Los Chankas"
$ws.Range("F23").Value = "W"
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 2
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 2.76
$ws.Range("L23").Value = 0.23
$ws.Range("M23").Value = 21
$ws.Range("N23").Value = 7
$ws.Range("O23").Value = 11
$ws.Range("P23").Value = 2
